# Edit script for "roles de metodologia scrum.docx"
# Implements:
#  1. New title paragraph "ROLES METODOLOGIA SCRUM" (bold, centered) inserted
#     at the very top, followed by a new centered/bold empty paragraph that
#     now hosts the "_GoBack" bookmark (moved up from the end of the doc),
#     followed by the (now single-run) "1. Dueño del producto (Product Owner)"
#     heading.
#  2. Merge split runs (same visible text, just consolidated into one run)
#     in a handful of paragraphs.
#  3. Drop the stray <w:lastRenderedPageBreak/> before "Relación de roles ".
#  4. Remove the old "_GoBack" bookmark pair near the end of the document
#     (it now lives at the top, see step 1).

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-OpenXmlPackage($bodyXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document ' + $wNs + '><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- Step 1: insert two new paragraphs before the current first paragraph ---
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphBefore()
$firstPara.Range.InsertParagraphBefore()

# Paragraph 1: "ROLES METODOLOGIA SCRUM" - bold, centered
$titleBody = '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr>' +
             '<w:r><w:rPr><w:b/></w:rPr><w:t>ROLES METODOLOGIA SCRUM</w:t></w:r></w:p>'
$d.Paragraphs(1).Range.InsertXML((New-OpenXmlPackage $titleBody))

# Paragraph 2: empty, bold+centered pPr, carries the relocated "_GoBack" bookmark
$bmBody = '<w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:b/></w:rPr></w:pPr>' +
          '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$d.Paragraphs(2).Range.InsertXML((New-OpenXmlPackage $bmBody))

# --- Step 2: consolidate the runs of "1. Due" + "ño del producto (Product Owner)" ---
$d.Content.Find.Execute("1. Dueño del producto (Product Owner)", $true, $false, $false, $false, $false, $true, 1, $false, "1. Dueño del producto (Product Owner)", 2) | Out-Null

# --- Step 3: merge the split run pairs (identical resulting text, single run) ---
$t1 = "Es el rol central del proyecto. En algunas ocasiones es quien representa al cliente y en otras son la misma persona. Sus principales funciones son:"
$d.Content.Find.Execute($t1, $true, $false, $false, $false, $false, $true, 1, $false, $t1, 2) | Out-Null

$t2 = "Muchos consideran que el Scrum Máster es el líder del proyecto. Y en cierta forma lo es, pero su labor no acaba ahí. Además de la gestión de las acciones en cada iteración, es el encargado de mantener en contacto al equipo de trabajo con el cliente. Otras de sus funciones más destacadas son:"
$d.Content.Find.Execute($t2, $true, $false, $false, $false, $false, $true, 1, $false, $t2, 2) | Out-Null

$t3 = "Finalmente, el Scrum Team hace referencia al equipo de trabajo que lleva a cabo las acciones propias de cada iteración: programadores, diseñadores, arquitectos, personal de servicio, entre otros. Lo principal es que deben estar organizados como un grupo o equipo y desempeñar roles concretos dentro de él. Se ocupan básicamente de cosas como las siguientes:"
$d.Content.Find.Execute($t3, $true, $false, $false, $false, $false, $true, 1, $false, $t3, 2) | Out-Null

# --- Step 4: drop the stale <w:lastRenderedPageBreak/> before "Relación de roles " ---
# Re-typing the paragraph's own text forces the run (and its cached
# lastRenderedPageBreak marker) to be rebuilt without that marker, while
# keeping the paragraph/run bold formatting intact.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -eq "Relación de roles ") {
        $p.Range.Text = "Relación de roles "
        break
    }
}

# --- Step 5: remove the old "_GoBack" bookmark paragraph near the end ---
# (its pPr/paragraph mark stays; only the bookmark pair is dropped - the
# bookmark now lives at the top of the document, see Step 1). That paragraph
# is the (normally empty) one right after the second bibliography hyperlink.
$xmlNoBookmark = New-OpenXmlPackage '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:rPr><w:b/></w:rPr></w:pPr></w:p>'
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*obs-edu*") {
        $d.Paragraphs($i + 1).Range.InsertXML($xmlNoBookmark)
        break
    }
}

Write-Output "done"
